$d = $word.ActiveDocument

# 1. The first paragraph holds the placeholder id text (with a trailing
#    space run after it). Replace the id text + trailing space in one pass
#    with the new id text - this both renames the placeholder and removes
#    the now-unwanted trailing-space run.
$d.Content.Find.Execute("**ID__AFFARS_5312_topic_11__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5312_70__ID**", 2)

# 2. Give that same paragraph a thin paragraph border and bump its left
#    indent from 6pt (120 twips) to 11.25pt (225 twips).
$para1 = $d.Paragraphs(1)
$pf = $para1.Format
$pf.LeftIndent = 11.25

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
